$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.644.89"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.57%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.316.81"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.56%  "

$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.44%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.71%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.634"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.61%  "

$ws.Range("E8").Value = "  -0.11%  "

$ws.Range("E9").Value = "  +1.20%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.09"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.13%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0913"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.63%  "

$ws.Range("E12").Value = "  +0.12%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.107"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.24%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.976"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.38%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.46"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.06%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.664.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.64%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.307.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.44%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.505.90"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.33%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.09%  "

$ws.Range("E20").Value = "  +0.94%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.66"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.85%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.72"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.17%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "280.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.55%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +20.32%  "

$ws.Range("E25").Value = "  -0.14%  "

$ws.Range("E26").Value = "  -0.17%  "

$ws.Range("E27").Value = "  -1.88%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.38"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.68%  "

$ws.Range("E29").Value = "  -0.36%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.64"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.09%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "165.57"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.67%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0881"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.71%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.43%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.137"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.98%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.71"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.84%  "

$ws.Range("E36").Value = "  -2.88%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.68"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.73%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0353"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.37%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.74"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.85%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "99.74"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.38%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.50"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.68%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "69.81"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.30%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.227"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.65%  "

$ws.Range("E45").Value = "  -0.20%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.04%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "113.78"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.52%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "78.87"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.50%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.03"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.37%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.33"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.04%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.615.79"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.40%  "
